$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = 45875
$ws.Range("A31").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("B31").Value = 6345.06005859375
$ws.Range("C31").Value = 6352.830078125
$ws.Range("D31").Value = 6301.10986328125
$ws.Range("E31").Value = 6309.2998046875
$ws.Range("F31").Value = 5408560000
$ws.Range("G31").Value = 0.0072819072950924
$ws.Range("H31").Value = 1

$ws.Range("I31").Value = 'when independent contractor are actually employee lesson from the m steadfast medical staffing judgment coursera cour q earnings call transcript ap business summarybrief at pm edt cohu q revenue rise wayfair post m q profit revenue up to b trump narrow fed chair candidate to four excluding treasury secretary bessent trump narrow fed chair candidate to four excluding treasury secretary bessent service sector slows down in july rural grocery store pilot project aim to prevent food desert in north dakota rural grocery store pilot project aim to prevent food desert in north dakota texas community q net jump trump say bank discriminate against conservativesincluding him the wall street journal microsoft cloud and ai revenue soar tevogenai expands collaboration with databricks and microsoft to build the beta version of it predictcelltm model with a focus on oncology trump administration officially end thing email for federal employee cbs news citizen bet on open banking to ease account switching best free tech book to read online cftc eye spot crypto trading on future market public opinion due august arabica coffee price sharply higher on brazil coffee crop concern danica mason named among top female entrepreneur build equity from the ground up with red team go franklin resource q revenue beat fincen seeing scam urge bank beware crypto atm mel gordon appointed ceo of chronicle heritage signaling a new era of stability and scalable growth cardano ada face new challenger a rollblock presale heat up with x speculation weave q revenue jump ensign ensg q earnings call transcript real estate expert kim nicols offer insight on buying luxury real estate in granite bay in hellonation magazine former x ceo yaccarino take helm at glp focused telehealth firm emed reuters california marijuana grower said employee were among hundred detained in immigration raid california marijuana grower said employee were among hundred detained in immigration raid u pressure tsmc to buy intel stake for taiwan tariff relief vertical aerospace move beyond can it fly with aciturri deal in focus linda yaccarino join health tech platform emed a ceo after leaving x techcrunch intels a chip woe low yield delay panther lake production gold technical setup turn bullish with job miss and cpi ahead michael novogratz say bitcoin ethereum treasury frenzy may be over what now join sdccu for free resource on national financial awareness day investor alert investigation of bellring brand inc brbr announced by holzer holzer llc gainey mckenna egleston announces a class action lawsuit ha been filed against lineage bed bath beyond could be headed your way here whats next for the revived retail chain bed bath beyond could be headed your way here whats next for the revived retail chain real good food announces corporate action requesting transfer of share to transfer agent johnson fistel investigates fairness of proposed sale of staar surgical johnson fistel commences investigation of agilon health inc birdie bringing froyo back with plan for location in west village week of july morning news rating today move back to first in demo ethereum news today ethereum etf see record m outflow a blackrock lead withdrawal gainey mckenna egleston announces a class action lawsuit ha been filed against lineage inc line sprott inc declares second quarter dividend featuring locus control report find landfill methane monitoring and tuning system among lowest cost ghg reduction sprott inc declares second quarter dividend middleton co inc ma sell share of procter gamble company the nysepg thoroughbred financial service llc acquires share of procter gamble company the nysepg palantir ceo warns of america ai danger zone a he plan to bring superpower to bluecollar worker fortune roku expands streaming option with new adfree service at month century complete expands to nevada with new home in pahrump replenish nutrient announces beiseker facility update crypto news today sharplink boost eth holding figure plan public listing galaxy mull tokenized stock paramount skydance board to include david ellison sherry lansing oracle safra catz the fda is targeting orange juice here why and what the science say the washington post introducing transparent lab grassfed protein shake the protein shake thor drink built to scale how barak diskin bridge real estate and startup growth westfield capital management co lp decrease stock position in alphabet inc nasdaqgoog alphabet inc nasdaqgoog position cut by atrium wealth solution inc syntegra private wealth group llc ha position in alphabet inc nasdaqgoog inspire medical system analyst lower their forecast after q result meta illegally collected data from flo period and pregnancy app jury find ar technica examining the future niagen bioscience earnings outlook earnings outlook for jack in the box a look ahead innovative ind props earnings forecast insight into helmerich paynes upcoming earnings a look ahead horace mann educatorss earnings forecast whats next honest co earnings preview what are the best smart glass of herbalife earnings preview rosen law firm urge novo nordisk a nyse nvo investor with loss in excess of k to contact the firm for information about their right a glimpse of genco shipping trading earnings potential a look ahead fastlys earnings forecast examining the future national vision holdingss earnings outlook earnings outlook for oscar health spotlight on mercadolibre analyzing the surge in option activity this is what whale are betting on philip morris intl a look into citigroup inc price over earnings is the market bullish or bearish on msc industrial direct co how is the market feeling about alphabet infoblox appoints joshua husk a chief revenue officer to accelerate global growth whats driving the market sentiment around caretrust reit bullet blockchain announces strategic initiative to accelerate growth and cement market leadership wyoming average gas price fall by cent per gallon in past week shiba inu shib may miss the mark but analyst say ruvi ai ruvi coinmarketcap listing powered surge ha it closer than ever novo nordisk expands legal action to protect u patient from unsafe nonfdaapproved compounded semaglutide pr newswire infoblox appoints joshua husk a chief revenue officer to accelerate global growth bullet blockchain announces strategic initiative to accelerate growth and cement market leadership beer institute release june taxable removal estimate report porsche gt r with manthey package arrives at porsche track experience buzz wear off for michigan marijuana business best meme coin to invest in today little pepe lilpepe to rally like pepe coin pepe did in cold wallet cashback utility set new standard while sol aim for doge eye target william raveis named top luxury brokerage at inman golden i club award beer institute release june taxable removal estimate report trump administration end musk controversial five thing email for federal worker live now ncsl hold discussion on citizenship and voting ag completes expansion of cold case unit trump give coal a boost gut renewables and will return haze to the four corner opinion trump narrow fed chair candidate to four excluding treasury secretary bessent track live update michigan primary election result track live result detroit mayoral primary election endoflease eviction allowed under new state law editorial quire served community with honor compassion texas ag asks judge to officially remove awol dems from seat this texas republican stalked trump to get his endorsementand flopped trump firing of bureau of labor statistic commissioner will prove damaging woof watch don lemon bonkers trumptrashing therapy session with the runaway texas dem disaster video federal judge blast republican for trying to pas obviously unconstitutional law city council approves limited shortterm rental in centertown overlay sjc confines zoning standing analysis to actual proposed use speculation a to future us is irrelevant time for yall to love america again texas dems in illinois call on republican to be brave ap news summary at pm edt trump take an unexpected walk on the white house roof to survey new project trump plan to take over gaza aid effort u official say opinion with harris out who lead california governor race pentagon keep a lid on golden dome marjorie taylor greene call for george santos year sentence to be commuted editorial the fed wa right to say no on interest rate banihashemi what possibility is there for true reform in iran garrity boiled frog syndrome they decided to cheat national democrat gov jb pritzker vow to do whatever it take to win redistricting battle well well what do we have here ethic committee rip aoc for spouse shenanigan and hooboy thread ayotte sign bill to tighten absentee voting ensure accessible voting exportsmouth teen pleads guilty to charge related to hatemotivated vandalism spree trump admin unveils rule to unleash american drone dominance eagle running back saquon barkley decline invitation to join trump sport council tulsi gabbard tell pod force one she belief there may be alien and still ha a lot of question on nj drone trump will restore a confederate monument in the nation capital justice department to seek federal hate crime charge and death penalty in killing of israeli embassy staffer the state of america judiciary with judge j michael luttig there it is jasmine crockett admits dems lose when voter id is the law of the land watch election is on democrat see opportunity and a wide open primary why spanish fighter jet are patrolling iceland sky raising deposit insurance limit would cost consumer over billion henderson state lawmaker running for dina titus house seat ct put back on federal list of sanctuary jurisdiction why the state top lawyer disagrees house committee issue subpoena for epstein file and deposition with the clinton republican next redistricting target missouri port will kirsten baesler be confirmed for a spot in the trump administration trump jack smith probe isnt just hypocritical it could backfire spectacularly u government proposes easing some restriction on drone traveling long distance live now state department spokesperson hold press briefing aug marjorie taylor greene call for george santos year sentence to be commuted marjorie taylor greene call for george santos year sentence to be commuted trump creates white house task force for los angeles olympics chuck todd worry democratic party hasnt lost enough to learn lesson change course aws announces general availability of amazon elastic vmware service trump plan drug tariff of up to live now texas house hold legislative session day after democrat flee state gop congressman face heated town hall where hundred boo him for supporting trump big bill the latest house committee issue subpoena for epstein file deposition with clinton and others palantir book it first billion in quarterly sale and dodge doge axe elizabeth warren go to the mat for mamdani but it the new york post with the takedown eazy grease accelerates multiregional expansion with liquid recovery solution merger building dominant uco collection network african nation have agreed to take deportee from the u what we know about the secretive deal african nation have agreed to take deportee from the u what we know about the secretive deal african nation have agreed to take deportee from the u what we know about the secretive deal african nation have agreed to take deportee from the u what we know about the secretive deal african nation have agreed to take deportee from the u what we know about the secretive deal african nation have agreed to take deportee from the u what we know about the secretive deal african nation have agreed to take deportee from the u what we know about the secretive deal african nation have agreed to take deportee from the u what we know about the secretive deal a former rolling stone say the met ha his stolen guitar the museum dispute it sterlington welcome leading private wealth team from morgan lewis new snap rule candy and soda no longer covered in additional state trump and democrat turn rival texas and california into proxy for national power struggle tariff trade and tension donald trump central asia policy massachusetts police training material classify mom for liberty a hate group alongside antifa the trump administration dismisses most on a federal board overseeing puerto rico finance the trump administration dismisses most on a federal board overseeing puerto rico finance the trump administration dismisses most on a federal board overseeing puerto rico finance the trump administration dismisses most on a federal board overseeing puerto rico finance the trump administration dismisses most on a federal board overseeing puerto rico finance the trump administration dismisses most on a federal board overseeing puerto rico finance the trump administration dismisses most on a federal board overseeing puerto rico finance the trump administration dismisses most on a federal board overseeing puerto rico finance the trump administration dismisses most on a federal board overseeing puerto rico finance aug might be one of the shortest day of the year here why from the publisher how to gain enemy and sow mistrust trump take an unexpected walk on the white house roof to survey new project trump take an unexpected walk on the white house roof to survey new project mayor brandon johnson say police fire pension sweetener will make chicago finance far more difficult staa alert monsey firm of wohl fruchter investigating fairness of the sale of staar surgical to alcon justice department open grand jury probe into investigation of trump and russia report nato to coordinate regular and largescale arm delivery to ukraine most will be bought in the u alligator alcatraz ha some new competition the speedway slammer pbgh announces aurora chen a vice president of advisory service sen michael bennet asks fed to allow colorado plan that would let toddler stay on medicaid tesla billion muskcentric gamble marjorie taylor greene urge trump to commute george santos sentence illinois gov pritzker commends texas democrat who fled state to block gop redistricting plan illinois gov pritzker commends texas democrat who fled state to block gop redistricting plan trump tour white house roof hint at further renovation woke men mocked by woman for obnoxious performative drinking habit you guy ruined matcha'
